$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new receipt rows for hours worked (use raw date serial numbers,
# matching the existing rows above)
$ws.Range("A63").Value = 40256
$ws.Range("B63").Value = 2
$ws.Range("C63").Value = "OMAP UART documentation"

$ws.Range("A64").Value = 40259
$ws.Range("B64").Value = 3
$ws.Range("C64").Value = "QNX 3D implementation design"

# Match date style/format of the existing date column
$ws.Range("A62").Copy()
$ws.Range("A63:A64").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A65").Select()
$excel.ActiveWindow.ScrollRow = 37
